$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds text dates formatted DD-MM-YYYY (e.g. "04-10-2021"). Assigning
# such a string straight to Range.Value gets auto-parsed as a real date by the
# COM layer (ambiguous MM-DD-YYYY read), which would store a date serial instead
# of literal text. Route it through a scratch formula cell (text result, not
# auto-typed) and PasteSpecial xlPasteValues (-4163) so the literal string lands
# in the target cell untouched, matching the existing shared-string cells above it.

$ws.Range("Z1").Formula = "=""04-10-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A191").PasteSpecial(-4163)
$ws.Range("B191").Value = 330.2
$ws.Range("C191").Value = 405
$ws.Range("D191").Value = 206
$ws.Range("E191").Value = 326
$ws.Range("F191").Value = 86.6
$ws.Range("G191").Value = 118.8
$ws.Range("H191").Value = 15.4
$ws.Range("I191").Value = 162.4
$ws.Range("J191").Value = 511.5
$ws.Range("K191").Value = 1612.6
$ws.Range("L191").Value = 310.3
$ws.Range("M191").Value = 301
$ws.Range("N191").Value = 156
$ws.Range("O191").Value = 364
$ws.Range("P191").Value = 186

$ws.Range("Z1").Formula = "=""05-10-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A192").PasteSpecial(-4163)
$ws.Range("B192").Value = 328.4
$ws.Range("C192").Value = 400
$ws.Range("D192").Value = 208
$ws.Range("E192").Value = 325
$ws.Range("F192").Value = 87
$ws.Range("G192").Value = 119.5
$ws.Range("H192").Value = 12.3
$ws.Range("I192").Value = 162.5
$ws.Range("J192").Value = 508.4
$ws.Range("K192").Value = 1591.6
$ws.Range("L192").Value = 309
$ws.Range("M192").Value = 297
$ws.Range("N192").Value = 153
$ws.Range("O192").Value = 358
$ws.Range("P192").Value = 183

$ws.Range("Z1").Formula = "=""06-10-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A193").PasteSpecial(-4163)
$ws.Range("B193").Value = 333
$ws.Range("C193").Value = 402
$ws.Range("D193").Value = 217
$ws.Range("E193").Value = 328
$ws.Range("F193").Value = 87.7
$ws.Range("G193").Value = 121.6
$ws.Range("H193").Value = 12.1
$ws.Range("I193").Value = 165
$ws.Range("J193").Value = 511.4
$ws.Range("K193").Value = 1588.1
$ws.Range("L193").Value = 310
$ws.Range("M193").Value = 298
$ws.Range("N193").Value = 156
$ws.Range("O193").Value = 363
$ws.Range("P193").Value = 185

$ws.Range("Z1").Formula = "=""07-10-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A194").PasteSpecial(-4163)
$ws.Range("B194").Value = 325.9
$ws.Range("C194").Value = 396
$ws.Range("D194").Value = 211
$ws.Range("E194").Value = 318
$ws.Range("F194").Value = 87.7
$ws.Range("G194").Value = 120.9
$ws.Range("H194").Value = 10.6
$ws.Range("I194").Value = 158.1
$ws.Range("J194").Value = 499.1
$ws.Range("K194").Value = 1582.1
$ws.Range("L194").Value = 303.1
$ws.Range("M194").Value = 289
$ws.Range("N194").Value = 152
$ws.Range("O194").Value = 356
$ws.Range("P194").Value = 173

$ws.Range("Z1").Formula = "=""08-10-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A195").PasteSpecial(-4163)
$ws.Range("B195").Value = 326.4
$ws.Range("C195").Value = 396
$ws.Range("D195").Value = 212
$ws.Range("E195").Value = 319
$ws.Range("F195").Value = 88.1
$ws.Range("G195").Value = 121.3
$ws.Range("H195").Value = 10.2
$ws.Range("I195").Value = 157.1
$ws.Range("J195").Value = 505.4
$ws.Range("K195").Value = 1602.4
$ws.Range("L195").Value = 302.3
$ws.Range("M195").Value = 286
$ws.Range("N195").Value = 153
$ws.Range("O195").Value = 356
$ws.Range("P195").Value = 171

$ws.Range("Z1").ClearContents()
$excel.CutCopyMode = 0
